$wb = $excel.ActiveWorkbook

# --- Restricciones_del_lider ---
$ws = $wb.Worksheets.Item("Restricciones_del_lider")

$r = $ws.Range("A2"); $r.NumberFormat = "@"; $r.Value = "2.09 - x"; $r.Style = "Normal"
$r = $ws.Range("B2"); $r.NumberFormat = "@"; $r.Value = "-3.09"; $r.Style = "Normal"
$r = $ws.Range("D2"); $r.NumberFormat = "@"; $r.Value = "0.86"; $r.Style = "Normal"

$r = $ws.Range("A3"); $r.NumberFormat = "@"; $r.Value = "-2.09 + x"; $r.Style = "Normal"
$r = $ws.Range("B3"); $r.NumberFormat = "@"; $r.Value = "1.0899999999999999"; $r.Style = "Normal"
$r = $ws.Range("D3"); $r.NumberFormat = "@"; $r.Value = "0.62"; $r.Style = "Normal"

$r = $ws.Range("A4"); $r.NumberFormat = "@"; $r.Value = "41.02289999999999 + x - y - 9(x^2)"; $r.Style = "Normal"
$r = $ws.Range("B4"); $r.NumberFormat = "@"; $r.Value = "-40.02289999999999"; $r.Style = "Normal"
$r = $ws.Range("D4"); $r.NumberFormat = "@"; $r.Value = "0.58"; $r.Style = "Normal"

# --- Restricciones_del_follower ---
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

$r = $ws.Range("A2"); $r.NumberFormat = "@"; $r.Value = "-22.9596 + (-0.5 + x)*(y^2)"; $r.Style = "Normal"
$r = $ws.Range("B2"); $r.NumberFormat = "@"; $r.Value = "22.9596"; $r.Style = "Normal"
$r = $ws.Range("D2"); $r.NumberFormat = "@"; $r.Value = "0.69"; $r.Style = "Normal"
$r = $ws.Range("E2"); $r.NumberFormat = "@"; $r.Value = "8.100000000000001"; $r.Style = "Normal"
$r = $ws.Range("F2"); $r.NumberFormat = "@"; $r.Value = "5.699999999999999"; $r.Style = "Normal"

$r = $ws.Range("A3"); $r.NumberFormat = "@"; $r.Value = "-3.8 + y"; $r.Style = "Normal"
$r = $ws.Range("B3"); $r.NumberFormat = "@"; $r.Value = "2.8"; $r.Style = "Normal"
$r = $ws.Range("D3"); $r.NumberFormat = "@"; $r.Value = "0.65"; $r.Style = "Normal"
$r = $ws.Range("E3"); $r.NumberFormat = "@"; $r.Value = "7.3"; $r.Style = "Normal"
$r = $ws.Range("F3"); $r.NumberFormat = "@"; $r.Value = "5.4"; $r.Style = "Normal"

$r = $ws.Range("A4"); $r.NumberFormat = "@"; $r.Value = "-5.8 - y"; $r.Style = "Normal"
$r = $ws.Range("B4"); $r.NumberFormat = "@"; $r.Value = "-4.8"; $r.Style = "Normal"
$r = $ws.Range("D4"); $r.NumberFormat = "@"; $r.Value = "0.32"; $r.Style = "Normal"
$r = $ws.Range("E4"); $r.NumberFormat = "@"; $r.Value = "7.0"; $r.Style = "Normal"
$r = $ws.Range("F4"); $r.NumberFormat = "@"; $r.Value = "9.8"; $r.Style = "Normal"

# --- Punto_modificado ---
$ws = $wb.Worksheets.Item("Punto_modificado")
$r = $ws.Range("A2"); $r.NumberFormat = "@"; $r.Value = "2.09"; $r.Style = "Normal"
$r = $ws.Range("B2"); $r.NumberFormat = "@"; $r.Value = "3.8"; $r.Style = "Normal"

# --- Vector_bf ---
# Note: worksheet lookup by name is case-insensitive, and "Vector_bf" /
# "Vector_BF" differ only by case, so we must address them positionally
# (sheet 5 = Vector_bf, sheet 6 = Vector_BF) to avoid ambiguity.
$ws = $wb.Worksheets.Item(5)
$r = $ws.Range("A2"); $r.NumberFormat = "@"; $r.Value = "-9.667959999999997"; $r.Style = "Normal"

# --- Vector_BF ---
$ws = $wb.Worksheets.Item(6)
$r = $ws.Range("A2"); $r.NumberFormat = "@"; $r.Value = "-96.48440000000002"; $r.Style = "Normal"
$r = $ws.Range("A3"); $r.NumberFormat = "@"; $r.Value = "-97.6004"; $r.Style = "Normal"
